$wb = $excel.ActiveWorkbook

# --- Rename sheets "01-02 R1" / "01-02 R2" to add " NO TRIP" suffix ---
# Renaming changes the sheet's internal name; re-assert the PageSetup.PrintArea
# afterwards so the workbook-level Print_Area defined name picks up the new
# sheet name (it otherwise keeps pointing at the stale name).
$ws0102R1 = $wb.Worksheets.Item("01-02 R1")
$ws0102R1.Name = "01-02 R1 NO TRIP"
$ws0102R1.PageSetup.PrintArea = '$A$1:$V$44'

$ws0102R2 = $wb.Worksheets.Item("01-02 R2")
$ws0102R2.Name = "01-02 R2 NO TRIP"
$ws0102R2.PageSetup.PrintArea = '$A$1:$V$44'

# --- Move the active/selected tab from "01-10 R3" to "01-02 R3" ---
$wsActive = $wb.Worksheets.Item("01-02 R3")
$wsActive.Activate()

# --- Data entry on "01-02 R3" ---
$ws = $wb.Worksheets.Item("01-02 R3")

# Row 8
$ws.Range("C8").Value = 1
$ws.Range("D8").Value = 4
$ws.Range("F8").Value = 1
$ws.Range("H8").Value = 18
$ws.Range("T8").Value = 2
$ws.Range("V8").Value = 6

# Row 15
$ws.Range("C15").Value = 1
$ws.Range("D15").Value = 4
$ws.Range("F15").Value = 0
$ws.Range("H15").Value = "13/6B"
$ws.Range("T15").Value = 2
$ws.Range("V15").Value = 0

# Row 16
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("H16").Value = "4/18B"
$ws.Range("T16").Value = 0
$ws.Range("V16").Value = 6

# Row 24
$ws.Range("H24").Value = 3
$ws.Range("J24").Value = 137
$ws.Range("L24").Value = 2
$ws.Range("M24").Value = 20
$ws.Range("N24").Value = 482
$ws.Range("Q24").Value = 20
$ws.Range("U24").Value = 8

# Row 31
$ws.Range("H31").Value = 3
$ws.Range("J31").Value = 62
$ws.Range("L31").Value = 2
$ws.Range("M31").Value = 8
$ws.Range("N31").Value = 2
$ws.Range("Q31").Value = 20
$ws.Range("U31").Value = 0

# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 75
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = 12
$ws.Range("N32").Value = 480
$ws.Range("Q32").Value = 0
$ws.Range("U32").Value = 8

# Row 36
$ws.Range("C36").Value = 1
$ws.Range("F36").Value = 680
$ws.Range("G36").Value = 2
$ws.Range("H36").Value = 21
$ws.Range("I36").Value = 6
$ws.Range("J36").Value = 6
$ws.Range("K36").Value = 14
$ws.Range("L36").Value = 4
$ws.Range("M36").Value = 1
$ws.Range("N36").Value = 6
$ws.Range("O36").Value = 3
$ws.Range("P36").Value = 2
$ws.Range("Q36").Value = 24

# Row 42
$ws.Range("C42").Value = 1
$ws.Range("F42").Value = 680
$ws.Range("G42").Value = 2
$ws.Range("H42").Value = 21
$ws.Range("I42").Value = 6
$ws.Range("J42").Value = 6
$ws.Range("K42").Value = 14
$ws.Range("L42").Value = 4
$ws.Range("M42").Value = 1
$ws.Range("N42").Value = 6
$ws.Range("O42").Value = 3
$ws.Range("P42").Value = 2
$ws.Range("Q42").Value = 24
